$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

# Update the repaymentstrategy value cell (B17) with the new scenario text
$ws.Range("B17").Value = "Overdue/Due Fee/Int,Principal"

# Move/refresh the active selection to B17 to match the saved view state
$ws.Activate()
$ws.Range("B17").Select()
